$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Table 1"
$ws.Range("B2").Value = "['Seat #1 is occupied by Izabela', 'Seat #2 is occupied by Maarten', 'Seat #3 is occupied by Petra', 'Seat #4 is occupied by Zelimkhan']"

$ws.Range("A3").Value = "Table 2"
$ws.Range("B3").Value = "['Seat #1 is occupied by Adheeba', 'Seat #2 is occupied by Kevin', 'Seat #3 is occupied by Soha', 'Seat #4 is occupied by Wouter']"

$ws.Range("A4").Value = "Table 3"
$ws.Range("B4").Value = "['Seat #1 is occupied by Kelli', 'Seat #2 is occupied by Minh Duc', 'Seat #3 is occupied by Nicolaas', 'Seat #4 is unoccupied']"

$ws.Range("A5").Value = "Table 4"
$ws.Range("B5").Value = "['Seat #1 is occupied by Anastasiia', 'Seat #2 is occupied by Muntadher', 'Seat #3 is occupied by Rasmita', 'Seat #4 is occupied by Yusra']"

$ws.Range("A6").Value = "Table 5"
$ws.Range("B6").Value = "['Seat #1 is occupied by Ihor', 'Seat #2 is occupied by Levin', 'Seat #3 is occupied by Tom', 'Seat #4 is occupied by Veena']"

$ws.Range("A7").Value = "Table 6"
$ws.Range("B7").Value = "['Seat #1 is occupied by Dhrisya', 'Seat #2 is occupied by Majid', 'Seat #3 is occupied by Rik', 'Seat #4 is occupied by Yeliz']"

$ws.Range("A8").Value = "Table 7"
$ws.Range("B8").Value = "['Seat #1 is occupied by Basma', 'Seat #2 is occupied by Moustafa', 'Seat #3 is occupied by Urson', 'Seat #4 is unoccupied']"
